# Add a "doc_ids" column to the "Tableau1" table on Sheet1, fill in the
# two known values, resize the new column, and restore the view state
# (frozen-pane scroll position + active selection) to match the saved
# workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Extend the table with a new last column -----------------------------
$newCol = $lo.ListColumns.Add()
$ws.Range("E1").Value = "doc_ids"

# Match the wrap-text formatting used by the rest of the table body.
$ws.Range("E1:E42").WrapText = $true

# --- Populate the new column's data ---------------------------------------
$ws.Range("E5").Value = "pdf_online"
$ws.Range("E8").Value = "pdf_online, bevnat_info"

# --- Column width -----------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 11.3

# --- Restore view state: scrolled-to row + active cell ---------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 2
$ws.Range("E8").Select() | Out-Null
